# Update lattice multiplication exercise table to new output at c986bee
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$v = [char]11  # vertical-tab = <w:br/> manual line break inside a Word Range.Text assignment

$c = $t.Cell(1, 1)
$c.Range.Text = "95 x 35" + $v + "  3    5" + $v + "  ----" + $v + "9|    |" + $v + "5|    |"

$c = $t.Cell(1, 2)
$c.Range.Text = "73 x 51" + $v + "  5    1" + $v + "  ----" + $v + "7|    |" + $v + "3|    |"

$c = $t.Cell(1, 3)
$c.Range.Text = "28 x 34" + $v + "  3    4" + $v + "  ----" + $v + "2|    |" + $v + "8|    |"

$c = $t.Cell(2, 1)
$c.Range.Text = "80 x 81" + $v + "  8    1" + $v + "  ----" + $v + "8|    |" + $v + "0|    |"

$c = $t.Cell(2, 2)
$c.Range.Text = "55 x 35" + $v + "  3    5" + $v + "  ----" + $v + "5|    |" + $v + "5|    |"

$c = $t.Cell(2, 3)
$c.Range.Text = "71 x 53" + $v + "  5    3" + $v + "  ----" + $v + "7|    |" + $v + "1|    |"

$c = $t.Cell(3, 1)
$c.Range.Text = "15 x 39" + $v + "  3    9" + $v + "  ----" + $v + "1|    |" + $v + "5|    |"

$c = $t.Cell(3, 2)
$c.Range.Text = "53 x 62" + $v + "  6    2" + $v + "  ----" + $v + "5|    |" + $v + "3|    |"

$c = $t.Cell(3, 3)
$c.Range.Text = "58 x 71" + $v + "  7    1" + $v + "  ----" + $v + "5|    |" + $v + "8|    |"

$c = $t.Cell(4, 1)
$c.Range.Text = "18 x 60" + $v + "  6    0" + $v + "  ----" + $v + "1|    |" + $v + "8|    |"

$c = $t.Cell(4, 2)
$c.Range.Text = "10 x 57" + $v + "  5    7" + $v + "  ----" + $v + "1|    |" + $v + "0|    |"

$c = $t.Cell(4, 3)
$c.Range.Text = "82 x 18" + $v + "  1    8" + $v + "  ----" + $v + "8|    |" + $v + "2|    |"

$c = $t.Cell(5, 1)
$c.Range.Text = "81 x 39" + $v + "  3    9" + $v + "  ----" + $v + "8|    |" + $v + "1|    |"

$c = $t.Cell(5, 2)
$c.Range.Text = "53 x 72" + $v + "  7    2" + $v + "  ----" + $v + "5|    |" + $v + "3|    |"

$c = $t.Cell(5, 3)
$c.Range.Text = "78 x 65" + $v + "  6    5" + $v + "  ----" + $v + "7|    |" + $v + "8|    |"

Write-Output "Updated 15 lattice multiplication cells."
